# Auto-generated Excel COM-interop script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.534.96"
$ws.Range("E2").Value = "  -0.67%  "

$ws.Range("D3").Value = "1.832.76"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  -0.09%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "312.77"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.26%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4292"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.67%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3667"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +0.40%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.07284"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -0.73%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.8654"
$cell.ClearFormats()
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.67"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.948.22"
$ws.Range("E12").Value = "  +5.58%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.397"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +0.59%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.544"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.12%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.06937"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -0.23%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.04%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "80.68"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +1.04%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.000008899"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -1.33%  "

$ws.Range("E19").Value = "  +0.02%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "15.41"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("D21").Value = "27.971.16"
$ws.Range("E21").Value = "  +0.98%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "5.147"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +3.25%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.82"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +4.59%  "

$ws.Range("D24").Value = "2.171.09"
$ws.Range("E24").Value = "  +2.99%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.995"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "154.42"
$cell.ClearFormats()
$ws.Range("E26").Value = "  -0.91%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "18.88"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +1.37%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "5.110"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -2.85%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "114.42"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -4.43%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.826"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -2.88%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08851"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.7505"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -1.09%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.986"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +0.79%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.535"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.43%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.132"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.55%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -1.80%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.05324"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -2.43%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01935"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -0.24%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.798"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.1663"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.5071"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -0.47%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "6.538"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -1.38%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.314"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -1.07%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "10.36"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.15%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "105.79"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +0.06%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.06482"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -1.13%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.4680"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.55%  "

$ws.Range("E49").Value = "  -0.03%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "1.615"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -1.75%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "63.72"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -1.48%  "
